$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 117.625
$ws.Range("I33").Value = 68.5
$ws.Range("J33").Value = 265
$ws.Range("K33").Value = 68.5
$ws.Range("L33").Value = 265
$ws.Range("M33").Value = 160.5
$ws.Range("N33").Value = -723
$ws.Range("H64").Value = 629897.5600000001
$ws.Range("I64").Value = 1003502.1
$ws.Range("J64").Value = 7223.3335
$ws.Range("K64").Value = 1003502.1
$ws.Range("L64").Value = 7223.3335
$ws.Range("M64").Value = -1003254.1
$ws.Range("N64").Value = -7719.3335
$ws.Range("H67").Value = 629897.5600000001
$ws.Range("I67").Value = 1003502.1
$ws.Range("J67").Value = 7223.3335
$ws.Range("K67").Value = 1003502.1
$ws.Range("L67").Value = 7223.3335
$ws.Range("M67").Value = -1002644.1
$ws.Range("N67").Value = -8939.333500000001
$ws.Range("H113").Value = 131313.88
$ws.Range("I113").Value = 174300.83
$ws.Range("J113").Value = 2353
$ws.Range("K113").Value = 174300.83
$ws.Range("L113").Value = 2353
$ws.Range("M113").Value = -171046.83
$ws.Range("N113").Value = -8861
$ws.Range("H133").Value = 13332
$ws.Range("J133").Value = 13332
$ws.Range("L133").Value = 13332
$ws.Range("N133").Value = -23452

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1013.375
$ws.Range("I2").Value = 691.82355
$ws.Range("J2").Value = 1794.2858
$ws.Range("K2").Value = 691.82355
$ws.Range("L2").Value = 1794.2858
$ws.Range("M2").Value = -578.82355
$ws.Range("N2").Value = -2020.2858
$ws.Range("H32").Value = 4112.727
$ws.Range("I32").Value = 2149.535
$ws.Range("K32").Value = 2149.535
$ws.Range("M32").Value = -1862.535
$ws.Range("H116").Value = 1013.375
$ws.Range("I116").Value = 691.82355
$ws.Range("J116").Value = 1794.2858
$ws.Range("K116").Value = 691.82355
$ws.Range("L116").Value = 1794.2858
$ws.Range("M116").Value = 1602.17645
$ws.Range("N116").Value = -6382.2858
$ws.Range("H122").Value = 2065.2856
$ws.Range("I122").Value = 1724.8
$ws.Range("K122").Value = 5174.4
$ws.Range("M122").Value = -2724.4
$ws.Range("H133").Value = 34400
$ws.Range("J133").Value = 34400
$ws.Range("L133").Value = 34400
$ws.Range("N133").Value = -39460
$ws.Range("H139").Value = 61243
$ws.Range("J139").Value = 61243
$ws.Range("L139").Value = 61243
$ws.Range("N139").Value = -71523

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1013.375
$ws.Range("I3").Value = 691.82355
$ws.Range("J3").Value = 1794.2858
$ws.Range("K3").Value = 691.82355
$ws.Range("L3").Value = 1794.2858
$ws.Range("M3").Value = -577.82355
$ws.Range("N3").Value = -2022.2858
$ws.Range("H22").Value = 278.8125
$ws.Range("I22").Value = 274.3846
$ws.Range("K22").Value = 274.3846
$ws.Range("M22").Value = -101.3846
$ws.Range("H86").Value = 10253.772
$ws.Range("I86").Value = 1842.1111
$ws.Range("J86").Value = 16077.23
$ws.Range("K86").Value = 1842.1111
$ws.Range("L86").Value = 16077.23
$ws.Range("M86").Value = -719.1111000000001
$ws.Range("N86").Value = -18323.23
$ws.Range("H89").Value = 10253.772
$ws.Range("I89").Value = 1842.1111
$ws.Range("J89").Value = 16077.23
$ws.Range("K89").Value = 9210.5555
$ws.Range("L89").Value = 80386.14999999999
$ws.Range("M89").Value = -3594.5555
$ws.Range("N89").Value = -91618.14999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 80.818184
$ws.Range("I7").Value = 58.42857
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 58.42857
$ws.Range("L7").Value = 120
$ws.Range("M7").Value = 54.57143
$ws.Range("N7").Value = -346
$ws.Range("H16").Value = 2011
$ws.Range("I16").Value = 2011
$ws.Range("K16").Value = 2011
$ws.Range("M16").Value = -1724
$ws.Range("H32").Value = 14558.167
$ws.Range("I32").Value = 14558.167
$ws.Range("K32").Value = 14558.167
$ws.Range("M32").Value = -14242.167
$ws.Range("H58").Value = 3152.52
$ws.Range("I58").Value = 2391.5557
$ws.Range("K58").Value = 2391.5557
$ws.Range("M58").Value = -2188.5557
$ws.Range("H62").Value = 15716.556
$ws.Range("I62").Value = 19741.584
$ws.Range("J62").Value = 7666.5
$ws.Range("K62").Value = 19741.584
$ws.Range("L62").Value = 7666.5
$ws.Range("M62").Value = -19117.584
$ws.Range("N62").Value = -8914.5
$ws.Range("H65").Value = 15716.556
$ws.Range("I65").Value = 19741.584
$ws.Range("J65").Value = 7666.5
$ws.Range("K65").Value = 98707.92
$ws.Range("L65").Value = 38332.5
$ws.Range("M65").Value = -95587.92
$ws.Range("N65").Value = -44572.5
$ws.Range("H99").Value = 4808932.5
$ws.Range("I99").Value = 6945680
$ws.Range("J99").Value = 1250
$ws.Range("K99").Value = 6945680
$ws.Range("L99").Value = 1250
$ws.Range("M99").Value = -6944182
$ws.Range("N99").Value = -4246
$ws.Range("H105").Value = 510.625
$ws.Range("I105").Value = 508.5
$ws.Range("J105").Value = 514.1667
$ws.Range("K105").Value = 508.5
$ws.Range("L105").Value = 514.1667
$ws.Range("M105").Value = 1238.5
$ws.Range("N105").Value = -4008.1667
$ws.Range("H107").Value = 410.7647
$ws.Range("I107").Value = 186.66667
$ws.Range("K107").Value = 186.66667
$ws.Range("M107").Value = 1733.33333
$ws.Range("H113").Value = 2011
$ws.Range("I113").Value = 2011
$ws.Range("K113").Value = 2011
$ws.Range("M113").Value = 159
$ws.Range("H122").Value = 1921.5294
$ws.Range("J122").Value = 3202.5
$ws.Range("L122").Value = 9607.5
$ws.Range("N122").Value = -14507.5
$ws.Range("H126").Value = 4808932.5
$ws.Range("I126").Value = 6945680
$ws.Range("J126").Value = 1250
$ws.Range("K126").Value = 20837040
$ws.Range("L126").Value = 3750
$ws.Range("M126").Value = -20834570
$ws.Range("N126").Value = -8690
$ws.Range("H132").Value = 3834.476
$ws.Range("I132").Value = 3074.182
$ws.Range("J132").Value = 4670.8
$ws.Range("K132").Value = 9222.545999999998
$ws.Range("L132").Value = 14012.4
$ws.Range("M132").Value = -6692.545999999998
$ws.Range("N132").Value = -19072.4
$ws.Range("H135").Value = 24375
$ws.Range("J135").Value = 24375
$ws.Range("L135").Value = 24375
$ws.Range("N135").Value = -34515
$ws.Range("H136").Value = 3152.52
$ws.Range("I136").Value = 2391.5557
$ws.Range("K136").Value = 7174.6671
$ws.Range("M136").Value = -4624.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 6760
$ws.Range("J102").Value = 6760
$ws.Range("L102").Value = 20280
$ws.Range("N102").Value = -25148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1027.5834
$ws.Range("I97").Value = 1116.5
$ws.Range("J97").Value = 938.6667
$ws.Range("K97").Value = 1116.5
$ws.Range("L97").Value = 938.6667
$ws.Range("M97").Value = -620.5
$ws.Range("N97").Value = -1930.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5793.6294
$ws.Range("I61").Value = 6087.8237
$ws.Range("J61").Value = 5293.5
$ws.Range("K61").Value = 6087.8237
$ws.Range("L61").Value = 5293.5
$ws.Range("M61").Value = -5885.8237
$ws.Range("N61").Value = -5697.5
$ws.Range("H93").Value = 1071.5
$ws.Range("I93").Value = 876.5
$ws.Range("J93").Value = 1364
$ws.Range("K93").Value = 876.5
$ws.Range("L93").Value = 1364
$ws.Range("M93").Value = 371.5
$ws.Range("N93").Value = -3860
$ws.Range("H100").Value = 2250.111
$ws.Range("I100").Value = 1684.2106
$ws.Range("J100").Value = 2882.5881
$ws.Range("K100").Value = 1684.2106
$ws.Range("L100").Value = 2882.5881
$ws.Range("M100").Value = -1143.2106
$ws.Range("N100").Value = -3964.5881
$ws.Range("H113").Value = 5793.6294
$ws.Range("I113").Value = 6087.8237
$ws.Range("J113").Value = 5293.5
$ws.Range("K113").Value = 6087.8237
$ws.Range("L113").Value = 5293.5
$ws.Range("M113").Value = -3917.8237
$ws.Range("N113").Value = -9633.5
$ws.Range("H139").Value = 53333.332
$ws.Range("J139").Value = 53333.332
$ws.Range("L139").Value = 53333.332
$ws.Range("N139").Value = -63613.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 2561
$ws.Range("I21").Value = 1833.3334
$ws.Range("J21").Value = 2997.6
$ws.Range("K21").Value = 1833.3334
$ws.Range("L21").Value = 2997.6
$ws.Range("M21").Value = -1598.3334
$ws.Range("N21").Value = -3467.6
$ws.Range("H35").Value = 2561
$ws.Range("I35").Value = 1833.3334
$ws.Range("J35").Value = 2997.6
$ws.Range("K35").Value = 1833.3334
$ws.Range("L35").Value = 2997.6
$ws.Range("M35").Value = -1543.3334
$ws.Range("N35").Value = -3577.6
$ws.Range("H113").Value = 412.9
$ws.Range("I113").Value = 194.18182
$ws.Range("J113").Value = 1014.375
$ws.Range("K113").Value = 582.5454599999999
$ws.Range("L113").Value = 3043.125
$ws.Range("M113").Value = 1587.45454
$ws.Range("N113").Value = -7383.125
